$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1545
$ws.Range("I11").Value = 1545
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1545
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1405

$ws.Range("H28").Value = 1957.9412
$ws.Range("I28").Value = 1184.3334
$ws.Range("J28").Value = 3814.6
$ws.Range("K28").Value = 1184.3334
$ws.Range("L28").Value = 3814.6
$ws.Range("M28").Value = -699.3334
$ws.Range("N28").Value = -4784.6

$ws.Range("H40").Value = 11113281
$ws.Range("I40").Value = 2115.8333
$ws.Range("J40").Value = 27780028
$ws.Range("K40").Value = 2115.8333
$ws.Range("L40").Value = 27780028
$ws.Range("M40").Value = -1940.8333
$ws.Range("N40").Value = -27780378

$ws.Range("H92").Value = 705.8570999999999
$ws.Range("I92").Value = 434.6111
$ws.Range("J92").Value = 2333.3333
$ws.Range("K92").Value = 434.6111
$ws.Range("L92").Value = 2333.3333
$ws.Range("M92").Value = 813.3888999999999
$ws.Range("N92").Value = -4829.3333

$ws.Range("H112").Value = 29758.658
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 29758.658
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 89275.974
$ws.Range("N112").Value = -91491.974
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1032.6154
$ws.Range("I2").Value = 1068.875
$ws.Range("J2").Value = 974.6
$ws.Range("K2").Value = 1068.875
$ws.Range("L2").Value = 974.6
$ws.Range("M2").Value = -955.875
$ws.Range("N2").Value = -1200.6

$ws.Range("H32").Value = 2904.9023
$ws.Range("I32").Value = 1879.2285
$ws.Range("J32").Value = 8888
$ws.Range("K32").Value = 1879.2285
$ws.Range("L32").Value = 8888
$ws.Range("M32").Value = -1592.2285
$ws.Range("N32").Value = -9462

$ws.Range("H45").Value = 912897.4
$ws.Range("I45").Value = 1252422.2
$ws.Range("J45").Value = 7497.6665
$ws.Range("K45").Value = 1252422.2
$ws.Range("L45").Value = 7497.6665
$ws.Range("M45").Value = -1252045.2
$ws.Range("N45").Value = -8251.666499999999

$ws.Range("H88").Value = 1237.6923
$ws.Range("I88").Value = 972.2857
$ws.Range("J88").Value = 1547.3334
$ws.Range("K88").Value = 972.2857
$ws.Range("L88").Value = 1547.3334
$ws.Range("M88").Value = -566.2857
$ws.Range("N88").Value = -2359.3334

$ws.Range("H91").Value = 1237.6923
$ws.Range("I91").Value = 972.2857
$ws.Range("J91").Value = 1547.3334
$ws.Range("K91").Value = 972.2857
$ws.Range("L91").Value = 1547.3334
$ws.Range("M91").Value = 431.7143
$ws.Range("N91").Value = -4355.3334

$ws.Range("H116").Value = 1032.6154
$ws.Range("I116").Value = 1068.875
$ws.Range("J116").Value = 974.6
$ws.Range("K116").Value = 1068.875
$ws.Range("L116").Value = 974.6
$ws.Range("M116").Value = 1225.125
$ws.Range("N116").Value = -5562.6

$ws.Range("H132").Value = 2781543
$ws.Range("I132").Value = 3706911
$ws.Range("J132").Value = 5438.6665
$ws.Range("K132").Value = 11120733
$ws.Range("L132").Value = 16315.9995
$ws.Range("M132").Value = -11118203
$ws.Range("N132").Value = -21375.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1032.6154
$ws.Range("I3").Value = 1068.875
$ws.Range("J3").Value = 974.6
$ws.Range("K3").Value = 1068.875
$ws.Range("L3").Value = 974.6
$ws.Range("M3").Value = -954.875
$ws.Range("N3").Value = -1202.6

$ws.Range("H20").Value = 2373.3333
$ws.Range("I20").Value = 1890.125
$ws.Range("J20").Value = 3919.6
$ws.Range("K20").Value = 1890.125
$ws.Range("L20").Value = 3919.6
$ws.Range("M20").Value = -1643.125
$ws.Range("N20").Value = -4413.6

$ws.Range("H86").Value = 2791.0715
$ws.Range("I86").Value = 2836.5386
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 2836.5386
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -1713.5386
$ws.Range("N86").Value = -4446

$ws.Range("H89").Value = 2791.0715
$ws.Range("I89").Value = 2836.5386
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 14182.693
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -8566.692999999999
$ws.Range("N89").Value = -22232

$ws.Range("H99").Value = 2178.9
$ws.Range("I99").Value = 1994
$ws.Range("J99").Value = 2363.8
$ws.Range("K99").Value = 1994
$ws.Range("L99").Value = 2363.8
$ws.Range("M99").Value = -496
$ws.Range("N99").Value = -5359.8

$ws.Range("H107").Value = 201979.4
$ws.Range("I107").Value = 949.5
$ws.Range("J107").Value = 335999.34
$ws.Range("K107").Value = 949.5
$ws.Range("L107").Value = 335999.34
$ws.Range("M107").Value = 970.5
$ws.Range("N107").Value = -339839.34

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1000000000
$ws.Range("I132").Value = 1000000000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000000000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2999997470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 45830.09
$ws.Range("I5").Value = 47893.43
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 143680.29
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -143568.29
$ws.Range("N5").Value = -7724

$ws.Range("H120").Value = 28999
$ws.Range("I120").Value = 28499
$ws.Range("J120").Value = 29999
$ws.Range("K120").Value = 85497
$ws.Range("L120").Value = 89997
$ws.Range("M120").Value = -80659
$ws.Range("N120").Value = -99673

$ws.Range("H131").Value = 2554.111
$ws.Range("I131").Value = 1997.8334
$ws.Range("J131").Value = 3666.6667
$ws.Range("K131").Value = 5993.5002
$ws.Range("L131").Value = 11000.0001
$ws.Range("M131").Value = -953.5002000000004
$ws.Range("N131").Value = -21080.0001

$ws.Range("H135").Value = 45830.09
$ws.Range("I135").Value = 47893.43
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 431040.87
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -428505.87
$ws.Range("N135").Value = -27570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1678999
$ws.Range("I7").Value = 3344999.2
$ws.Range("J7").Value = 12998.667
$ws.Range("K7").Value = 3344999.2
$ws.Range("L7").Value = 12998.667
$ws.Range("M7").Value = -3344887.2
$ws.Range("N7").Value = -13222.667

$ws.Range("H8").Value = 1678999
$ws.Range("I8").Value = 3344999.2
$ws.Range("J8").Value = 12998.667
$ws.Range("K8").Value = 3344999.2
$ws.Range("L8").Value = 12998.667
$ws.Range("M8").Value = -3344860.2
$ws.Range("N8").Value = -13276.667

$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 30000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 30000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -29702
$ws.Range("N35").ClearContents()

$ws.Range("H40").Value = 35000
$ws.Range("I40").Value = 35000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 35000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -34849

$ws.Range("H44").Value = 9999
$ws.Range("I44").Value = 9999
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 9999
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -9403

$ws.Range("H97").Value = 740.1
$ws.Range("I97").Value = 464.82352
$ws.Range("J97").Value = 2300
$ws.Range("K97").Value = 464.82352
$ws.Range("L97").Value = 2300
$ws.Range("M97").Value = 31.17648000000003
$ws.Range("N97").Value = -3292

$ws.Range("H98").Value = 125460
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 125460
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 125460
$ws.Range("N98").Value = -131450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3957.7222
$ws.Range("I40").Value = 3827.4375
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3827.4375
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3691.4375
$ws.Range("N40").Value = -5272

$ws.Range("H43").Value = 138750
$ws.Range("I43").Value = 50000
$ws.Range("J43").Value = 151428.58
$ws.Range("K43").Value = 50000
$ws.Range("L43").Value = 151428.58
$ws.Range("M43").Value = -49807
$ws.Range("N43").Value = -151814.58

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
